$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("I2").Value = 0.5154172875118646
$ws.Range("J2").Value = 0.5154172875118646
$ws.Range("M2").Value = 1.321358333333333
$ws.Range("N2").Value = 3.964075
$ws.Range("O2").Value = 0.06904315418552966
$ws.Range("P2").Value = 0.06904315418552966
$ws.Range("Q2").Value = 0.08562490090555554
$ws.Range("R2").Value = 0.7706241081499999
$ws.Range("S2").Value = 0.03558603525156914
$ws.Range("T2").Value = 0.03558603525156914

$ws.Range("I3").Value = 0.5154172875118646
$ws.Range("J3").Value = 0.5154172875118646
$ws.Range("O3").Value = 0.4558096119837698
$ws.Range("P3").Value = 0.4558096119837698
$ws.Range("R3").Value = 5.08751200412
$ws.Range("S3").Value = 0.2349321538305101
$ws.Range("T3").Value = 0.2349321538305101

$ws.Range("I4").Value = 0.5154172875118646
$ws.Range("J4").Value = 0.5154172875118646
$ws.Range("M4").Value = 9.093439666666667
$ws.Range("N4").Value = 27.280319
$ws.Range("O4").Value = 0.4751472338307006
$ws.Range("P4").Value = 0.4751472338307005
$ws.Range("Q4").Value = 0.5892609526931111
$ws.Range("R4").Value = 5.303348574237999
$ws.Range("S4").Value = 0.2448990984297853
$ws.Range("T4").Value = 0.2448990984297853

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.060924
$ws.Range("H5").Value = 0.182772
$ws.Range("I5").Value = 0.4845827124881355
$ws.Range("J5").Value = 0.4845827124881354
$ws.Range("M5").Value = 1.321358333333333
$ws.Range("N5").Value = 3.964075
$ws.Range("O5").Value = 0.06904315418552966
$ws.Range("P5").Value = 0.06904315418552966
$ws.Range("Q5").Value = 0.08050243509999999
$ws.Range("R5").Value = 0.7245219158999999
$ws.Range("S5").Value = 0.03345711893396053
$ws.Range("T5").Value = 0.03345711893396052

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.060924
$ws.Range("H6").Value = 0.182772
$ws.Range("I6").Value = 0.4845827124881355
$ws.Range("J6").Value = 0.4845827124881354
$ws.Range("O6").Value = 0.4558096119837698
$ws.Range("P6").Value = 0.4558096119837698
$ws.Range("Q6").Value = 0.53146157848
$ws.Range("R6").Value = 4.78315420632
$ws.Range("S6").Value = 0.2208774581532597
$ws.Range("T6").Value = 0.2208774581532597

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.060924
$ws.Range("H7").Value = 0.182772
$ws.Range("I7").Value = 0.4845827124881355
$ws.Range("J7").Value = 0.4845827124881354
$ws.Range("M7").Value = 9.093439666666667
$ws.Range("N7").Value = 27.280319
$ws.Range("O7").Value = 0.4751472338307006
$ws.Range("P7").Value = 0.4751472338307005
$ws.Range("Q7").Value = 0.554008718252
$ws.Range("R7").Value = 4.986078464267999
$ws.Range("S7").Value = 0.2302481354009153
$ws.Range("T7").Value = 0.2302481354009152
